$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Text format on D2:D51 so numeric-looking strings
# (e.g. "1.013") are not auto-converted to numbers by Excel.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '29.002.84'
$ws.Range("D3").Value = '1.969.50'
$ws.Range("D5").Value = '329.20'
$ws.Range("D6").Value = '1.013'
$ws.Range("D7").Value = '0.4966'
$ws.Range("D8").Value = '0.4179'
$ws.Range("D9").Value = '54.34'
$ws.Range("D10").Value = '0.09270'
$ws.Range("D11").Value = '1.091'
$ws.Range("D12").Value = '22.65'
$ws.Range("D13").Value = '1.984.50'
$ws.Range("D14").Value = '7.869'
$ws.Range("D15").Value = '6.437'
$ws.Range("D16").Value = '1.015'
$ws.Range("D17").Value = '0.00001108'
$ws.Range("D18").Value = '91.46'
$ws.Range("D19").Value = '0.06739'
$ws.Range("D20").Value = '19.07'
$ws.Range("D21").Value = '1.013'
$ws.Range("D22").Value = '5.953'
$ws.Range("D23").Value = '29.024.76'
$ws.Range("D24").Value = '11.91'
$ws.Range("D25").Value = '2.269'
$ws.Range("D26").Value = '2.221.32'
$ws.Range("D27").Value = '20.70'
$ws.Range("D28").Value = '156.61'
$ws.Range("D29").Value = '6.221'
$ws.Range("D30").Value = '2.252'
$ws.Range("D31").Value = '127.07'
$ws.Range("D32").Value = '1.041'
$ws.Range("D33").Value = '0.09810'
$ws.Range("D34").Value = '1.501'
$ws.Range("D35").Value = '5.787'
$ws.Range("D36").Value = '3.745'
$ws.Range("D37").Value = '0.02412'
$ws.Range("D38").Value = '1.317'
$ws.Range("D39").Value = '0.06374'
$ws.Range("D40").Value = '9.001'
$ws.Range("D41").Value = '0.6448'
$ws.Range("D42").Value = '11.40'
$ws.Range("D43").Value = '0.1995'
$ws.Range("D44").Value = '1.013'
$ws.Range("D45").Value = '0.6164'
$ws.Range("D46").Value = '1.346'
$ws.Range("D47").Value = '13.32'
$ws.Range("D48").Value = '2.163'
$ws.Range("D49").Value = '3.485'
$ws.Range("D51").Value = '0.06944'

# Restore the original (default) cell style now that values are set,
# so no stray style/number-format is left behind on these cells.
$dRange.Style = "Normal"

# E column values already contain non-numeric padding (spaces, % sign)
# so they remain text automatically.
$ws.Range("E2").Value = '  -1.68%  '
$ws.Range("E3").Value = '  -1.26%  '
$ws.Range("E4").Value = '  +0.72%  '
$ws.Range("E5").Value = '  -0.15%  '
$ws.Range("E6").Value = '  +0.57%  '
$ws.Range("E7").Value = '  -0.67%  '
$ws.Range("E9").Value = '  +4.89%  '
$ws.Range("E10").Value = '  +4.24%  '
$ws.Range("E11").Value = '  -2.65%  '
$ws.Range("E12").Value = '  -2.96%  '
$ws.Range("E13").Value = '  -3.20%  '
$ws.Range("E14").Value = '  -2.96%  '
$ws.Range("E15").Value = '  -1.16%  '
$ws.Range("E16").Value = '  +0.72%  '
$ws.Range("E17").Value = '  +0.01%  '
$ws.Range("E18").Value = '  -4.97%  '
$ws.Range("E19").Value = '  +1.56%  '
$ws.Range("E20").Value = '  -3.38%  '
$ws.Range("E22").Value = '  -0.27%  '
$ws.Range("E23").Value = '  -1.71%  '
$ws.Range("E24").Value = '  +0.15%  '
$ws.Range("E25").Value = '  -0.31%  '
$ws.Range("E26").Value = '  -2.41%  '
$ws.Range("E27").Value = '  +0.42%  '
$ws.Range("E28").Value = '  -0.86%  '
$ws.Range("E29").Value = '  -5.28%  '
$ws.Range("E30").Value = '  -3.55%  '
$ws.Range("E31").Value = '  -0.72%  '
$ws.Range("E32").Value = '  -1.11%  '
$ws.Range("E33").Value = '  -1.35%  '
$ws.Range("E34").Value = '  -3.70%  '
$ws.Range("E35").Value = '  -0.95%  '
$ws.Range("E36").Value = '  -1.14%  '
$ws.Range("E37").Value = '  -2.02%  '
$ws.Range("E38").Value = '  +2.22%  '
$ws.Range("E39").Value = '  +0.26%  '
$ws.Range("E40").Value = '  -6.17%  '
$ws.Range("E41").Value = '  -1.30%  '
$ws.Range("E42").Value = '  -2.80%  '
$ws.Range("E43").Value = '  -3.49%  '
$ws.Range("E44").Value = '  +0.72%  '
$ws.Range("E45").Value = '  -2.94%  '
$ws.Range("E46").Value = '  +5.83%  '
$ws.Range("E47").Value = '  -1.15%  '
$ws.Range("E48").Value = '  -2.34%  '
$ws.Range("E49").Value = '  -1.12%  '
$ws.Range("E50").Value = '  -1.20%  '
$ws.Range("E51").Value = '  -0.94%  '
